# C5-PowerPoint.pptx edit:
#  1. Slide 6: the "Sources of finance" table switches from the default
#     "Table_0" table style to PowerPoint's built-in "Medium Style 2 -
#     Accent 1" gallery style (tableStyleId GUID change).
#  2. The deck's theme colour scheme (ppt/theme/theme2.xml, the theme
#     actually bound to the slide master / slides) switches from the
#     custom "Integral" palette to the stock Office palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{E33ECDDE-8F0E-4B5F-BBF4-6907E93A34A4}")

# --- 2. Theme colours -------------------------------------------------
# ThemeColor.RGB uses the standard OLE/VBA BGR-packed integer, so to
# land a final a:srgbClr of RRGGBB we feed 0xBBGGRR.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0x000000   # dk1
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
